# Update backend to work with more stock tickers
# - Insert a "Company" column after Ticker
# - Add Warner Bros. Discovery (WBD) and Moderna (MRNA) rows, replacing
#   the previous Netflix (NFLX) and Microsoft (MSFT) rows
# - Re-point hyperlinks at the new column layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink objects first; we'll recreate them once
# the new column layout/values are in place (their cell refs shift).
$ws.Hyperlinks.Delete()

# Insert a new column B ("Company"); old B/C slide right to C/D.
$ws.Columns.Item(2).Insert()

# --- Header row ---
$ws.Range("B1").Value = "Company"

# --- Row 2 : GOOG / Alphabet Inc. ---
$ws.Range("B2").Value = "Alphabet Inc."

# --- Row 3 : META / Meta Platforms, Inc. ---
$ws.Range("B3").Value = "Meta Platforms, Inc."

# --- Row 4 : TEAM / Atlassian Corporation (wrap text) ---
$ws.Range("B4").Value = "Atlassian Corporation"
$ws.Range("B4").WrapText = $true

# --- Row 5 : was NFLX -> now WBD / Warner Bros. Discovery, Inc. ---
$ws.Range("A5").Value = "WBD"
$ws.Range("B5").Value = "Warner Bros. Discovery, Inc."
$ws.Range("C5").Value = "https://ir.corporate.discovery.com/investor-relations/default.aspx"
$ws.Range("D5").Value = "https://s201.q4cdn.com/336605034/files/doc_earnings/2023/q4/earnings-result/WBD-4Q23-Earnings-Release.pdf"
$ws.Range("D5").Style = "Normal"

# --- Row 6 : was MSFT -> now MRNA / Moderna, Inc. ---
$ws.Range("A6").Value = "MRNA"
$ws.Range("B6").Value = "Moderna, Inc."
$ws.Range("C6").Value = "https://investors.modernatx.com/overview/default.aspx"
$ws.Range("D6").Value = "https://s29.q4cdn.com/435878511/files/doc_earnings/2023/q4/earnings-result/Q4-23-PR_Final.pdf"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 31
$ws.Columns.Item(3).ColumnWidth = 59.8
$ws.Columns.Item(4).ColumnWidth = 103.25

# --- Re-create hyperlinks on the new layout ---
$ws.Hyperlinks.Add($ws.Range("C2"), "https://abc.xyz/investor/")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://abc.xyz/assets/95/eb/9cef90184e09bac553796896c633/2023q4-alphabet-earnings-release.pdf")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://investor.fb.com/home/default.aspx")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://s21.q4cdn.com/399680738/files/doc_financials/2023/q4/Meta-12-31-2023-Exhibit-99-1-FINAL.pdf")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://investors.atlassian.com/")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://s28.q4cdn.com/541786762/files/doc_financials/2024/q2/TEAM-Q2-2024-Earnings-Release.pdf")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://ir.corporate.discovery.com/investor-relations/default.aspx")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://investors.modernatx.com/overview/default.aspx")
$ws.Hyperlinks.Add($ws.Range("D6"), "https://s29.q4cdn.com/435878511/files/doc_earnings/2023/q4/earnings-result/Q4-23-PR_Final.pdf")
